$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.645.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.672.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4821"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2633"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06175"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.665.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6007"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.424"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9990"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9991"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.641.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006797"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.12%  "

$ws.Range("E20").Value = "  +1.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.476"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.877.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("E23").Value = "  +3.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.374"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("E26").Value = "  +3.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.406"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "105.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.99%  "

$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("E30").Value = "  +4.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.674"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07705"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04380"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9981"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6194"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9559"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.629"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8739"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01518"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.878"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.689"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1128"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.253"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05265"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.384"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3364"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "

